$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New opportunity listings scraped for this run (replaces rows 2-10; old row 11 is dropped)
$data = @(
    @('1326965', 'https://aiesec.org/opportunity/global-talent/1326965', 'Sales Specialist', 'Ankara, Türkiye', 'No', '62 applicants', '6 - 18 Months', 'Enveks Toz Kontrol Teknolojileri Sanayi ve Ticaret A.Ş.'),
    @('1326423', 'https://aiesec.org/opportunity/global-talent/1326423', 'Sales Specialist', 'Sincan, İstasyon, 06934 Sincan/Ankara, Türkiye', 'No', '40 applicants', '9 - 12 Weeks', 'JUMBOSAN MAKİNA İNŞAAT OTOMOTİV NAKLİYE VE GIDA SANAYİ TİCAR'),
    @('1325832', 'https://aiesec.org/opportunity/global-talent/1325832', 'Inventory Management Intern', 'Makati City, Metro Manila, Philippines', 'No', '16 applicants', '3 - 6 Months', 'Consistent Frozen Solutions Corporation'),
    @('1325829', 'https://aiesec.org/opportunity/global-talent/1325829', 'Business Process Management Intern', 'Makati City, Metro Manila, Philippines', 'No', '13 applicants', '3 - 6 Months', 'Consistent Frozen Solutions Corporation'),
    @('1308313', 'https://aiesec.org/opportunity/global-talent/1308313', 'Mechanical/Electrical/Chemical/Environmental Engineer Intern', 'Ankara, Türkiye', 'No', '186 applicants', '6 - 18 Months', 'AKA SU TEKNOLOJİLERİ ANONİM ŞİRKETİ'),
    @('1304372', 'https://aiesec.org/opportunity/global-talent/1304372', 'Marketing Specialist', 'Ankara, Türkiye', 'No', '396 applicants', '9 - 12 Weeks', 'Ganioğulları LTD'),
    @('1303648', 'https://aiesec.org/opportunity/global-talent/1303648', 'Social Media Manager', 'Ahangama, Sri Lanka', 'No', '137 applicants', '3 - 6 Months', 'Surfing Wombats'),
    @('1294219', 'https://aiesec.org/opportunity/global-talent/1294219', 'Marketing Specialist', 'Altındağ, Türkiye', 'No', '103 applicants', '6 - 18 Months', 'Torksan RPM Polat Makina'),
    @('1293818', 'https://aiesec.org/opportunity/global-talent/1293818', 'Receptionist and Marketing  Intern', 'Chalkidiki, Greece', 'No', '181 applicants', '9 - 12 Weeks', 'Hotel Kriopigi'),
)

# Opportunity IDs are textual identifiers, not numbers -- keep column A as text
$ws.Range("A2:A10").NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}

# Dataset shrank from 10 listings to 9 -- drop the now-unused last row
$ws.Rows.Item(11).Delete()

# Column widths adjusted to fit the refreshed content
# (values nudged slightly below the integer target to compensate for the
# character-width -> pixel -> character-width rounding Excel performs on save,
# so the saved width ends up exactly at the intended integer)
$ws.Columns.Item(3).ColumnWidth = 62.15
$ws.Columns.Item(4).ColumnWidth = 48.15
$ws.Columns.Item(6).ColumnWidth = 16.15

